$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.121.02"
$ws.Range("E2").Value = "  +0.01%  "

$ws.Range("D3").Value = "2.757.94"
$ws.Range("E3").Value = "  -0.56%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "575.83"
$ws.Range("E5").Value = "  -1.93%  "

$ws.Range("D6").Value = "159.63"
$ws.Range("E6").Value = "  -1.25%  "

$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("E8").Value = "  -3.02%  "

$ws.Range("E9").Value = "  -3.68%  "

$ws.Range("D10").Value = "5.90"
$ws.Range("E10").Value = "  -13.12%  "

$ws.Range("E11").Value = "  +3.42%  "

$ws.Range("E12").Value = "  -3.42%  "

$ws.Range("D13").Value = "3.247.31"

$ws.Range("D14").Value = "27.00"
$ws.Range("E14").Value = "  -2.56%  "

$ws.Range("D15").Value = "63.707.65"
$ws.Range("E15").Value = "  -0.47%  "

$ws.Range("E16").Value = "  -5.64%  "

$ws.Range("D17").Value = "2.763.35"
$ws.Range("E17").Value = "  -0.59%  "

$ws.Range("E18").Value = "  -2.64%  "

$ws.Range("D19").Value = "4.82"
$ws.Range("E19").Value = "  -4.69%  "

$ws.Range("D20").Value = "359.83"
$ws.Range("E20").Value = "  -2.09%  "

$ws.Range("D21").Value = "6.63"
$ws.Range("E21").Value = "  -6.30%  "

$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  -0.32%  "

$ws.Range("E23").Value = "  -7.72%  "

$ws.Range("D24").Value = "65.14"
$ws.Range("E24").Value = "  -3.82%  "

$ws.Range("E25").Value = "  -3.71%  "

$ws.Range("E26").Value = "  -2.89%  "

$ws.Range("E27").Value = "  -0.12%  "

$ws.Range("D29").Value = "7.38"
$ws.Range("E29").Value = "  +1.32%  "

$ws.Range("E30").Value = "  -4.18%  "

$ws.Range("D31").Value = "1.34"
$ws.Range("E31").Value = "  +4.83%  "

$ws.Range("D32").Value = "170.46"
$ws.Range("E32").Value = "  -1.17%  "

$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "4.94"
$ws.Range("E33").Value = "  -4.02%  "

$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "20.22"
$ws.Range("E34").Value = "  -3.04%  "

$ws.Range("E35").Value = "  -1.85%  "

$ws.Range("E36").Value = "  +0.09%  "

$ws.Range("E37").Value = "  -1.87%  "

$ws.Range("E38").Value = "  -2.24%  "

$ws.Range("D39").Value = "353.34"
$ws.Range("E39").Value = "  +3.14%  "

$ws.Range("D40").Value = "6.33"
$ws.Range("E40").Value = "  +0.31%  "

$ws.Range("E41").Value = "  -2.67%  "

$ws.Range("E42").Value = "  -2.82%  "

$ws.Range("E43").Value = "  -3.84%  "

$ws.Range("D44").Value = "21.99"
$ws.Range("E44").Value = "  -3.10%  "

$ws.Range("D45").Value = "0.0587"
$ws.Range("E45").Value = "  -4.10%  "

$ws.Range("D46").Value = "137.10"
$ws.Range("E46").Value = "  -1.15%  "

$ws.Range("E47").Value = "  -3.20%  "

$ws.Range("D48").Value = "0.628"
$ws.Range("E48").Value = "  -3.57%  "

$ws.Range("E49").Value = "  -1.86%  "

$ws.Range("D50").Value = "0.998"
$ws.Range("E50").Value = "  +0.03%  "

$ws.Range("E51").Value = "  +0.22%  "
